$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Cxcl12"
$ws.Cells.Item(2,3).Value = "Cxcr3"
$ws.Cells.Item(2,4).Value = "M1"
$ws.Cells.Item(2,5).Value = 3
$ws.Cells.Item(2,6).Value = 1
$ws.Cells.Item(2,7).Value = 109.1118546666667
$ws.Cells.Item(2,8).Value = 327.335564
$ws.Cells.Item(2,9).Value = 0.3029068882986101
$ws.Cells.Item(2,10).Value = 0.3029068882986101
$ws.Cells.Item(2,11).Value = 2
$ws.Cells.Item(2,12).Value = 0.6666666666666666
$ws.Cells.Item(2,13).Value = 0.8058999999999999
$ws.Cells.Item(2,14).Value = 2.4177
$ws.Cells.Item(2,15).Value = 0.09451456273674329
$ws.Cells.Item(2,16).Value = 0.0945145627367433
$ws.Cells.Item(2,17).Value = 87.93324367586666
$ws.Cells.Item(2,18).Value = 791.3991930828
$ws.Cells.Item(2,19).Value = 0.02862911209749067
$ws.Cells.Item(2,20).Value = 0.02862911209749068

# Row 3
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Cxcl12"
$ws.Cells.Item(3,3).Value = "Cxcr3"
$ws.Cells.Item(3,4).Value = "M2"
$ws.Cells.Item(3,5).Value = 3
$ws.Cells.Item(3,6).Value = 1
$ws.Cells.Item(3,7).Value = 109.1118546666667
$ws.Cells.Item(3,8).Value = 327.335564
$ws.Cells.Item(3,9).Value = 0.3029068882986101
$ws.Cells.Item(3,10).Value = 0.3029068882986101
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 2.375759
$ws.Cells.Item(3,14).Value = 7.127277
$ws.Cells.Item(3,15).Value = 0.2786249200308754
$ws.Cells.Item(3,16).Value = 0.2786249200308755
$ws.Cells.Item(3,17).Value = 259.2234707310253
$ws.Cells.Item(3,18).Value = 2333.011236579228
$ws.Cells.Item(3,19).Value = 0.08439740752900154
$ws.Cells.Item(3,20).Value = 0.08439740752900156

# Row 4
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Cxcl12"
$ws.Cells.Item(4,3).Value = "Cxcr3"
$ws.Cells.Item(4,4).Value = "Neutro"
$ws.Cells.Item(4,5).Value = 3
$ws.Cells.Item(4,6).Value = 1
$ws.Cells.Item(4,7).Value = 109.1118546666667
$ws.Cells.Item(4,8).Value = 327.335564
$ws.Cells.Item(4,9).Value = 0.3029068882986101
$ws.Cells.Item(4,10).Value = 0.3029068882986101
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 5.345069333333334
$ws.Cells.Item(4,14).Value = 16.035208
$ws.Cells.Item(4,15).Value = 0.6268605172323812
$ws.Cells.Item(4,16).Value = 0.6268605172323812
$ws.Cells.Item(4,17).Value = 583.2104282819237
$ws.Cells.Item(4,18).Value = 5248.893854537313
$ws.Cells.Item(4,19).Value = 0.1898803686721178
$ws.Cells.Item(4,20).Value = 0.1898803686721178

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Cxcl12"
$ws.Cells.Item(5,3).Value = "Cxcr3"
$ws.Cells.Item(5,4).Value = "M1"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 227.6338753333333
$ws.Cells.Item(5,8).Value = 682.901626
$ws.Cells.Item(5,9).Value = 0.6319374650831437
$ws.Cells.Item(5,10).Value = 0.6319374650831437
$ws.Cells.Item(5,11).Value = 2
$ws.Cells.Item(5,12).Value = 0.6666666666666666
$ws.Cells.Item(5,13).Value = 0.8058999999999999
$ws.Cells.Item(5,14).Value = 2.4177
$ws.Cells.Item(5,15).Value = 0.09451456273674329
$ws.Cells.Item(5,16).Value = 0.0945145627367433
$ws.Cells.Item(5,17).Value = 183.4501401311333
$ws.Cells.Item(5,18).Value = 1651.0512611802
$ws.Cells.Item(5,19).Value = 0.05972729318929931
$ws.Cells.Item(5,20).Value = 0.05972729318929932

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Cxcl12"
$ws.Cells.Item(6,3).Value = "Cxcr3"
$ws.Cells.Item(6,4).Value = "M2"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 227.6338753333333
$ws.Cells.Item(6,8).Value = 682.901626
$ws.Cells.Item(6,9).Value = 0.6319374650831437
$ws.Cells.Item(6,10).Value = 0.6319374650831437
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 2.375759
$ws.Cells.Item(6,14).Value = 7.127277
$ws.Cells.Item(6,15).Value = 0.2786249200308754
$ws.Cells.Item(6,16).Value = 0.2786249200308755
$ws.Cells.Item(6,17).Value = 540.8032280280446
$ws.Cells.Item(6,18).Value = 4867.229052252402
$ws.Cells.Item(6,19).Value = 0.1760735256733051
$ws.Cells.Item(6,20).Value = 0.1760735256733051

# Row 7
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Cxcl12"
$ws.Cells.Item(7,3).Value = "Cxcr3"
$ws.Cells.Item(7,4).Value = "Neutro"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 227.6338753333333
$ws.Cells.Item(7,8).Value = 682.901626
$ws.Cells.Item(7,9).Value = 0.6319374650831437
$ws.Cells.Item(7,10).Value = 0.6319374650831437
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 5.345069333333334
$ws.Cells.Item(7,14).Value = 16.035208
$ws.Cells.Item(7,15).Value = 0.6268605172323812
$ws.Cells.Item(7,16).Value = 0.6268605172323812
$ws.Cells.Item(7,17).Value = 1216.718846272023
$ws.Cells.Item(7,18).Value = 10950.46961644821
$ws.Cells.Item(7,19).Value = 0.3961366462205393
$ws.Cells.Item(7,20).Value = 0.3961366462205393

# Row 8
$ws.Cells.Item(8,1).Value = "M1"
$ws.Cells.Item(8,2).Value = "Cxcl12"
$ws.Cells.Item(8,3).Value = "Cxcr3"
$ws.Cells.Item(8,4).Value = "M1"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 0.194568
$ws.Cells.Item(8,8).Value = 0.583704
$ws.Cells.Item(8,9).Value = 0.0005401428435299865
$ws.Cells.Item(8,10).Value = 0.0005401428435299865
$ws.Cells.Item(8,11).Value = 2
$ws.Cells.Item(8,12).Value = 0.6666666666666666
$ws.Cells.Item(8,13).Value = 0.8058999999999999
$ws.Cells.Item(8,14).Value = 2.4177
$ws.Cells.Item(8,15).Value = 0.09451456273674329
$ws.Cells.Item(8,16).Value = 0.0945145627367433
$ws.Cells.Item(8,17).Value = 0.1568023512
$ws.Cells.Item(8,18).Value = 1.4112211608
$ws.Cells.Item(8,19).Value = 0.00005105136467161782
$ws.Cells.Item(8,20).Value = 0.00005105136467161783

# Row 9
$ws.Cells.Item(9,1).Value = "M1"
$ws.Cells.Item(9,2).Value = "Cxcl12"
$ws.Cells.Item(9,3).Value = "Cxcr3"
$ws.Cells.Item(9,4).Value = "M2"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 0.194568
$ws.Cells.Item(9,8).Value = 0.583704
$ws.Cells.Item(9,9).Value = 0.0005401428435299865
$ws.Cells.Item(9,10).Value = 0.0005401428435299865
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 2.375759
$ws.Cells.Item(9,14).Value = 7.127277
$ws.Cells.Item(9,15).Value = 0.2786249200308754
$ws.Cells.Item(9,16).Value = 0.2786249200308755
$ws.Cells.Item(9,17).Value = 0.462246677112
$ws.Cells.Item(9,18).Value = 4.160220094008
$ws.Cells.Item(9,19).Value = 0.0001504972565837921
$ws.Cells.Item(9,20).Value = 0.0001504972565837922

# Row 10
$ws.Cells.Item(10,1).Value = "M1"
$ws.Cells.Item(10,2).Value = "Cxcl12"
$ws.Cells.Item(10,3).Value = "Cxcr3"
$ws.Cells.Item(10,4).Value = "Neutro"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 0.194568
$ws.Cells.Item(10,8).Value = 0.583704
$ws.Cells.Item(10,9).Value = 0.0005401428435299865
$ws.Cells.Item(10,10).Value = 0.0005401428435299865
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 5.345069333333334
$ws.Cells.Item(10,14).Value = 16.035208
$ws.Cells.Item(10,15).Value = 0.6268605172323812
$ws.Cells.Item(10,16).Value = 0.6268605172323812
$ws.Cells.Item(10,17).Value = 1.039979450048
$ws.Cells.Item(10,18).Value = 9.359815050432001
$ws.Cells.Item(10,19).Value = 0.0003385942222745765
$ws.Cells.Item(10,20).Value = 0.0003385942222745765

# Row 11
$ws.Cells.Item(11,1).Value = "M2"
$ws.Cells.Item(11,2).Value = "Cxcl12"
$ws.Cells.Item(11,3).Value = "Cxcr3"
$ws.Cells.Item(11,4).Value = "M1"
$ws.Cells.Item(11,5).Value = 3
$ws.Cells.Item(11,6).Value = 1
$ws.Cells.Item(11,7).Value = 0.4517016666666667
$ws.Cells.Item(11,8).Value = 1.355105
$ws.Cells.Item(11,9).Value = 0.001253975076377243
$ws.Cells.Item(11,10).Value = 0.001253975076377243
$ws.Cells.Item(11,11).Value = 2
$ws.Cells.Item(11,12).Value = 0.6666666666666666
$ws.Cells.Item(11,13).Value = 0.8058999999999999
$ws.Cells.Item(11,14).Value = 2.4177
$ws.Cells.Item(11,15).Value = 0.09451456273674329
$ws.Cells.Item(11,16).Value = 0.0945145627367433
$ws.Cells.Item(11,17).Value = 0.3640263731666666
$ws.Cells.Item(11,18).Value = 3.2762373585
$ws.Cells.Item(11,19).Value = 0.0001185189060265694
$ws.Cells.Item(11,20).Value = 0.0001185189060265694

# Row 12
$ws.Cells.Item(12,1).Value = "M2"
$ws.Cells.Item(12,2).Value = "Cxcl12"
$ws.Cells.Item(12,3).Value = "Cxcr3"
$ws.Cells.Item(12,4).Value = "M2"
$ws.Cells.Item(12,5).Value = 3
$ws.Cells.Item(12,6).Value = 1
$ws.Cells.Item(12,7).Value = 0.4517016666666667
$ws.Cells.Item(12,8).Value = 1.355105
$ws.Cells.Item(12,9).Value = 0.001253975076377243
$ws.Cells.Item(12,10).Value = 0.001253975076377243
$ws.Cells.Item(12,11).Value = 3
$ws.Cells.Item(12,12).Value = 1
$ws.Cells.Item(12,13).Value = 2.375759
$ws.Cells.Item(12,14).Value = 7.127277
$ws.Cells.Item(12,15).Value = 0.2786249200308754
$ws.Cells.Item(12,16).Value = 0.2786249200308755
$ws.Cells.Item(12,17).Value = 1.073134299898333
$ws.Cells.Item(12,18).Value = 9.658208699085
$ws.Cells.Item(12,19).Value = 0.0003493887053763202
$ws.Cells.Item(12,20).Value = 0.0003493887053763203

# Row 13
$ws.Cells.Item(13,1).Value = "M2"
$ws.Cells.Item(13,2).Value = "Cxcl12"
$ws.Cells.Item(13,3).Value = "Cxcr3"
$ws.Cells.Item(13,4).Value = "Neutro"
$ws.Cells.Item(13,5).Value = 3
$ws.Cells.Item(13,6).Value = 1
$ws.Cells.Item(13,7).Value = 0.4517016666666667
$ws.Cells.Item(13,8).Value = 1.355105
$ws.Cells.Item(13,9).Value = 0.001253975076377243
$ws.Cells.Item(13,10).Value = 0.001253975076377243
$ws.Cells.Item(13,11).Value = 3
$ws.Cells.Item(13,12).Value = 1
$ws.Cells.Item(13,13).Value = 5.345069333333334
$ws.Cells.Item(13,14).Value = 16.035208
$ws.Cells.Item(13,15).Value = 0.6268605172323812
$ws.Cells.Item(13,16).Value = 0.6268605172323812
$ws.Cells.Item(13,17).Value = 2.414376726315556
$ws.Cells.Item(13,18).Value = 21.72939053684
$ws.Cells.Item(13,19).Value = 0.0007860674649743534
$ws.Cells.Item(13,20).Value = 0.0007860674649743534

# Row 14
$ws.Cells.Item(14,1).Value = "Neutro"
$ws.Cells.Item(14,2).Value = "Cxcl12"
$ws.Cells.Item(14,3).Value = "Cxcr3"
$ws.Cells.Item(14,4).Value = "M1"
$ws.Cells.Item(14,5).Value = 3
$ws.Cells.Item(14,6).Value = 1
$ws.Cells.Item(14,7).Value = 0.7729569999999999
$ws.Cells.Item(14,8).Value = 2.318871
$ws.Cells.Item(14,9).Value = 0.002145816331084288
$ws.Cells.Item(14,10).Value = 0.002145816331084288
$ws.Cells.Item(14,11).Value = 2
$ws.Cells.Item(14,12).Value = 0.6666666666666666
$ws.Cells.Item(14,13).Value = 0.8058999999999999
$ws.Cells.Item(14,14).Value = 2.4177
$ws.Cells.Item(14,15).Value = 0.09451456273674329
$ws.Cells.Item(14,16).Value = 0.0945145627367433
$ws.Cells.Item(14,17).Value = 0.6229260462999999
$ws.Cells.Item(14,18).Value = 5.606334416699999
$ws.Cells.Item(14,19).Value = 0.0002028108922457942
$ws.Cells.Item(14,20).Value = 0.0002028108922457943

# Row 15
$ws.Cells.Item(15,1).Value = "Neutro"
$ws.Cells.Item(15,2).Value = "Cxcl12"
$ws.Cells.Item(15,3).Value = "Cxcr3"
$ws.Cells.Item(15,4).Value = "M2"
$ws.Cells.Item(15,5).Value = 3
$ws.Cells.Item(15,6).Value = 1
$ws.Cells.Item(15,7).Value = 0.7729569999999999
$ws.Cells.Item(15,8).Value = 2.318871
$ws.Cells.Item(15,9).Value = 0.002145816331084288
$ws.Cells.Item(15,10).Value = 0.002145816331084288
$ws.Cells.Item(15,11).Value = 3
$ws.Cells.Item(15,12).Value = 1
$ws.Cells.Item(15,13).Value = 2.375759
$ws.Cells.Item(15,14).Value = 7.127277
$ws.Cells.Item(15,15).Value = 0.2786249200308754
$ws.Cells.Item(15,16).Value = 0.2786249200308755
$ws.Cells.Item(15,17).Value = 1.836359549363
$ws.Cells.Item(15,18).Value = 16.527235944267
$ws.Cells.Item(15,19).Value = 0.0005978779036493062
$ws.Cells.Item(15,20).Value = 0.0005978779036493063

# Row 16
$ws.Cells.Item(16,1).Value = "Neutro"
$ws.Cells.Item(16,2).Value = "Cxcl12"
$ws.Cells.Item(16,3).Value = "Cxcr3"
$ws.Cells.Item(16,4).Value = "Neutro"
$ws.Cells.Item(16,5).Value = 3
$ws.Cells.Item(16,6).Value = 1
$ws.Cells.Item(16,7).Value = 0.7729569999999999
$ws.Cells.Item(16,8).Value = 2.318871
$ws.Cells.Item(16,9).Value = 0.002145816331084288
$ws.Cells.Item(16,10).Value = 0.002145816331084288
$ws.Cells.Item(16,11).Value = 3
$ws.Cells.Item(16,12).Value = 1
$ws.Cells.Item(16,13).Value = 5.345069333333334
$ws.Cells.Item(16,14).Value = 16.035208
$ws.Cells.Item(16,15).Value = 0.6268605172323812
$ws.Cells.Item(16,16).Value = 0.6268605172323812
$ws.Cells.Item(16,17).Value = 4.131508756685333
$ws.Cells.Item(16,18).Value = 37.18357881016799
$ws.Cells.Item(16,19).Value = 0.001345127535189187
$ws.Cells.Item(16,20).Value = 0.001345127535189187

# Row 17
$ws.Cells.Item(17,1).Value = "sCs"
$ws.Cells.Item(17,2).Value = "Cxcl12"
$ws.Cells.Item(17,3).Value = "Cxcr3"
$ws.Cells.Item(17,4).Value = "M1"
$ws.Cells.Item(17,5).Value = 3
$ws.Cells.Item(17,6).Value = 1
$ws.Cells.Item(17,7).Value = 22.05086833333333
$ws.Cells.Item(17,8).Value = 66.15260499999999
$ws.Cells.Item(17,9).Value = 0.06121571236725463
$ws.Cells.Item(17,10).Value = 0.06121571236725463
$ws.Cells.Item(17,11).Value = 2
$ws.Cells.Item(17,12).Value = 0.6666666666666666
$ws.Cells.Item(17,13).Value = 0.8058999999999999
$ws.Cells.Item(17,14).Value = 2.4177
$ws.Cells.Item(17,15).Value = 0.09451456273674329
$ws.Cells.Item(17,16).Value = 0.0945145627367433
$ws.Cells.Item(17,17).Value = 17.77079478983333
$ws.Cells.Item(17,18).Value = 159.9371531085
$ws.Cells.Item(17,19).Value = 0.00578577628700932
$ws.Cells.Item(17,20).Value = 0.005785776287009321

# Row 18
$ws.Cells.Item(18,1).Value = "sCs"
$ws.Cells.Item(18,2).Value = "Cxcl12"
$ws.Cells.Item(18,3).Value = "Cxcr3"
$ws.Cells.Item(18,4).Value = "M2"
$ws.Cells.Item(18,5).Value = 3
$ws.Cells.Item(18,6).Value = 1
$ws.Cells.Item(18,7).Value = 22.05086833333333
$ws.Cells.Item(18,8).Value = 66.15260499999999
$ws.Cells.Item(18,9).Value = 0.06121571236725463
$ws.Cells.Item(18,10).Value = 0.06121571236725463
$ws.Cells.Item(18,11).Value = 3
$ws.Cells.Item(18,12).Value = 1
$ws.Cells.Item(18,13).Value = 2.375759
$ws.Cells.Item(18,14).Value = 7.127277
$ws.Cells.Item(18,15).Value = 0.2786249200308754
$ws.Cells.Item(18,16).Value = 0.2786249200308755
$ws.Cells.Item(18,17).Value = 52.38754890073166
$ws.Cells.Item(18,18).Value = 471.487940106585
$ws.Cells.Item(18,19).Value = 0.01705622296295939
$ws.Cells.Item(18,20).Value = 0.0170562229629594

# Row 19
$ws.Cells.Item(19,1).Value = "sCs"
$ws.Cells.Item(19,2).Value = "Cxcl12"
$ws.Cells.Item(19,3).Value = "Cxcr3"
$ws.Cells.Item(19,4).Value = "Neutro"
$ws.Cells.Item(19,5).Value = 3
$ws.Cells.Item(19,6).Value = 1
$ws.Cells.Item(19,7).Value = 22.05086833333333
$ws.Cells.Item(19,8).Value = 66.15260499999999
$ws.Cells.Item(19,9).Value = 0.06121571236725463
$ws.Cells.Item(19,10).Value = 0.06121571236725463
$ws.Cells.Item(19,11).Value = 3
$ws.Cells.Item(19,12).Value = 1
$ws.Cells.Item(19,13).Value = 5.345069333333334
$ws.Cells.Item(19,14).Value = 16.035208
$ws.Cells.Item(19,15).Value = 0.6268605172323812
$ws.Cells.Item(19,16).Value = 0.6268605172323812
$ws.Cells.Item(19,17).Value = 117.8634201018711
$ws.Cells.Item(19,18).Value = 1060.77078091684
$ws.Cells.Item(19,19).Value = 0.03837371311728591
$ws.Cells.Item(19,20).Value = 0.03837371311728591
